# "shift update on frontend"
# Adds a new "shift" column (M) to the upload template sheet, with a
# sample value of "morning shift" in the example data row, and moves
# the sheet's active selection over to the newly added column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell (M1) + sample data cell (M2) for the "shift" column.
# The engine manages the shared-string table automatically, so these
# plain assignments are enough to create the new <si> entries.
$ws.Range("M1").Value = "shift"
$ws.Range("M2").Value = "morning shift"

# Give the new column the same width as the other data columns
# (closest representable ColumnWidth to the 14.57-character width
# used elsewhere in the sheet).
$ws.Columns.Item(13).ColumnWidth = 13.67

# Move the selection over to where the new column now sits.
$ws.Range("N6").Select()
